$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Stephen Curry"
$ws.Range("B2").Value = "PG,SG"
$ws.Range("C2").Value = "Golden State Warriors"
$ws.Range("A3").Value = "Klay Thompson"
$ws.Range("B3").Value = "SG,SF"
$ws.Range("C3").Value = "Dallas Mavericks"
$ws.Range("A4").Value = "Austin Reaves"
$ws.Range("B4").Value = "PG,SG"
$ws.Range("C4").Value = "Los Angeles Lakers"
$ws.Range("A5").Value = "Darius Garland"
$ws.Range("B5").Value = "PG"
$ws.Range("C5").Value = "Cleveland Cavaliers"
$ws.Range("A6").Value = "Franz Wagner"
$ws.Range("B6").Value = "SF,PF"
$ws.Range("C6").Value = "Orlando Magic"
$ws.Range("A7").Value = "Keegan Murray"
$ws.Range("B7").Value = "SF,PF"
$ws.Range("C7").Value = "Sacramento Kings"
$ws.Range("A8").Value = "Kevin Durant"
$ws.Range("B8").Value = "SF,PF"
$ws.Range("C8").Value = "Phoenix Suns"
$ws.Range("A9").Value = "Cameron Johnson"
$ws.Range("B9").Value = "SF,PF"
$ws.Range("C9").Value = "Brooklyn Nets"
$ws.Range("A10").Value = "Jarrett Allen"
$ws.Range("B10").Value = "C"
$ws.Range("C10").Value = "Cleveland Cavaliers"
$ws.Range("A11").Value = "Karl-Anthony Towns"
$ws.Range("B11").Value = "PF,C"
$ws.Range("C11").Value = "New York Knicks"
$ws.Range("A12").Value = "Trey Murphy III"
$ws.Range("B12").Value = "SF,PF"
$ws.Range("C12").Value = "New Orleans Pelicans"
$ws.Range("A13").Value = "Jalen Duren"
$ws.Range("B13").Value = "C"
$ws.Range("C13").Value = "Detroit Pistons"
$ws.Range("A14").Value = "Tyrese Maxey"
$ws.Range("B14").Value = "PG,SG"
$ws.Range("C14").Value = "Philadelphia 76ers"
$ws.Range("A15").Value = "Tyrese Haliburton"
$ws.Range("B15").Value = "PG,SG"
$ws.Range("C15").Value = "Indiana Pacers"
$ws.Range("A16").Value = "Nick Richards"
$ws.Range("B16").Value = "C"
$ws.Range("C16").Value = "Phoenix Suns"
$ws.Range("A17").Value = "OG Anunoby"
$ws.Range("B17").Value = "SF,PF"
$ws.Range("C17").Value = "New York Knicks"
$ws.Range("A18").Value = "Mark Williams"
$ws.Range("B18").Value = "C"
$ws.Range("C18").Value = "Charlotte Hornets"
